# Scheduled-runner refresh of the Fenrir market-board snapshot.
# Re-pulls current Universalis average prices (columns H-L) for the
# tracked Leve items on each job sheet and recomputes the dependent
# NQ/HQ profit columns (M/N) from the refreshed prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 25539
$ws.Range("J57").Value = 26073.75
$ws.Range("L57").Value = 78221.25
$ws.Range("N57").Value = -79219.25
$ws.Range("H62").Value = 66668196
$ws.Range("I62").Value = 76924376
$ws.Range("J62").Value = 3003
$ws.Range("K62").Value = 76924376
$ws.Range("L62").Value = 3003
$ws.Range("M62").Value = -76923752
$ws.Range("N62").Value = -4251
$ws.Range("H65").Value = 66668196
$ws.Range("I65").Value = 76924376
$ws.Range("J65").Value = 3003
$ws.Range("K65").Value = 384621880
$ws.Range("L65").Value = 15015
$ws.Range("M65").Value = -384618760
$ws.Range("N65").Value = -21255
$ws.Range("H112").Value = 5922.5
$ws.Range("I112").Value = 715.8333
$ws.Range("J112").Value = 7484.5
$ws.Range("K112").Value = 2147.4999
$ws.Range("L112").Value = 22453.5
$ws.Range("M112").Value = -1039.4999
$ws.Range("N112").Value = -24669.5
$ws.Range("H132").Value = 26244102
$ws.Range("I132").Value = 31378568
$ws.Range("J132").Value = 1273.7778
$ws.Range("K132").Value = 94135704
$ws.Range("L132").Value = 3821.3334
$ws.Range("M132").Value = -94133174
$ws.Range("N132").Value = -8881.3334
$ws.Range("H137").Value = 20652532
$ws.Range("I137").Value = 426690.44
$ws.Range("J137").Value = 47620320
$ws.Range("K137").Value = 1280071.32
$ws.Range("L137").Value = 142860960
$ws.Range("M137").Value = -1277521.32
$ws.Range("N137").Value = -142866060
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1168.1818
$ws.Range("I2").Value = 835.7143
$ws.Range("J2").Value = 1750
$ws.Range("K2").Value = 835.7143
$ws.Range("L2").Value = 1750
$ws.Range("M2").Value = -722.7143
$ws.Range("N2").Value = -1976
$ws.Range("H45").Value = 582996.9399999999
$ws.Range("I45").Value = 715314.4399999999
$ws.Range("J45").Value = 800
$ws.Range("K45").Value = 715314.4399999999
$ws.Range("L45").Value = 800
$ws.Range("M45").Value = -714937.4399999999
$ws.Range("N45").Value = -1554
$ws.Range("H61").Value = 8775007
$ws.Range("I61").Value = 9526893
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 9526893
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -9526681
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 677.1579
$ws.Range("I74").Value = 459.5
$ws.Range("J74").Value = 1838
$ws.Range("K74").Value = 459.5
$ws.Range("L74").Value = 1838
$ws.Range("M74").Value = 414.5
$ws.Range("N74").Value = -3586
$ws.Range("H77").Value = 677.1579
$ws.Range("I77").Value = 459.5
$ws.Range("J77").Value = 1838
$ws.Range("K77").Value = 2297.5
$ws.Range("L77").Value = 9190
$ws.Range("M77").Value = 2070.5
$ws.Range("N77").Value = -17926
$ws.Range("H116").Value = 1168.1818
$ws.Range("I116").Value = 835.7143
$ws.Range("J116").Value = 1750
$ws.Range("K116").Value = 835.7143
$ws.Range("L116").Value = 1750
$ws.Range("M116").Value = 1458.2857
$ws.Range("N116").Value = -6338
$ws.Range("H136").Value = 8775007
$ws.Range("I136").Value = 9526893
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 28580679
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -28578129
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1168.1818
$ws.Range("I3").Value = 835.7143
$ws.Range("J3").Value = 1750
$ws.Range("K3").Value = 835.7143
$ws.Range("L3").Value = 1750
$ws.Range("M3").Value = -721.7143
$ws.Range("N3").Value = -1978
$ws.Range("H134").Value = 25661258
$ws.Range("I134").Value = 31770658
$ws.Range("J134").Value = 1780
$ws.Range("K134").Value = 95311974
$ws.Range("L134").Value = 5340
$ws.Range("M134").Value = -95309439
$ws.Range("N134").Value = -10410
$ws.Range("H139").Value = 61259.332
$ws.Range("J139").Value = 61259.332
$ws.Range("L139").Value = 61259.332
$ws.Range("N139").Value = -71539.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19613484
$ws.Range("I31").Value = 23810570
$ws.Range("J31").Value = 27084.5
$ws.Range("K31").Value = 23810570
$ws.Range("L31").Value = 27084.5
$ws.Range("M31").Value = -23810275
$ws.Range("N31").Value = -27674.5
$ws.Range("H34").Value = 19613484
$ws.Range("I34").Value = 23810570
$ws.Range("J34").Value = 27084.5
$ws.Range("K34").Value = 23810570
$ws.Range("L34").Value = 27084.5
$ws.Range("M34").Value = -23810368
$ws.Range("N34").Value = -27488.5
$ws.Range("H58").Value = 3200113.2
$ws.Range("I58").Value = 3996894
$ws.Range("J58").Value = 12990.444
$ws.Range("K58").Value = 3996894
$ws.Range("L58").Value = 12990.444
$ws.Range("M58").Value = -3996691
$ws.Range("N58").Value = -13396.444
$ws.Range("H99").Value = 83335380
$ws.Range("I99").Value = 121214470
$ws.Range("J99").Value = 1360
$ws.Range("K99").Value = 121214470
$ws.Range("L99").Value = 1360
$ws.Range("M99").Value = -121212972
$ws.Range("N99").Value = -4356
$ws.Range("H126").Value = 83335380
$ws.Range("I126").Value = 121214470
$ws.Range("J126").Value = 1360
$ws.Range("K126").Value = 363643410
$ws.Range("L126").Value = 4080
$ws.Range("M126").Value = -363640940
$ws.Range("N126").Value = -9020
$ws.Range("H132").Value = 9528620
$ws.Range("I132").Value = 13891011
$ws.Range("J132").Value = 10676.637
$ws.Range("K132").Value = 41673033
$ws.Range("L132").Value = 32029.911
$ws.Range("M132").Value = -41670503
$ws.Range("N132").Value = -37089.911
$ws.Range("H134").Value = 38604530
$ws.Range("I134").Value = 48078270
$ws.Range("J134").Value = 7814850
$ws.Range("K134").Value = 144234810
$ws.Range("L134").Value = 23444550
$ws.Range("M134").Value = -144232275
$ws.Range("N134").Value = -23449620
$ws.Range("H136").Value = 3200113.2
$ws.Range("I136").Value = 3996894
$ws.Range("J136").Value = 12990.444
$ws.Range("K136").Value = 11990682
$ws.Range("L136").Value = 38971.33199999999
$ws.Range("M136").Value = -11988132
$ws.Range("N136").Value = -44071.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1679.1428
$ws.Range("I129").Value = 487.14285
$ws.Range("J129").Value = 2871.1428
$ws.Range("K129").Value = 1461.42855
$ws.Range("L129").Value = 8613.428400000001
$ws.Range("M129").Value = 3538.57145
$ws.Range("N129").Value = -18613.4284
$ws.Range("H130").Value = 8335
$ws.Range("J130").Value = 1363.6364
$ws.Range("L130").Value = 4090.9092
$ws.Range("N130").Value = -14130.9092
$ws.Range("H131").Value = 15260361
$ws.Range("I131").Value = 52631970
$ws.Range("J131").Value = 1605350.6
$ws.Range("K131").Value = 157895910
$ws.Range("L131").Value = 4816051.800000001
$ws.Range("M131").Value = -157890870
$ws.Range("N131").Value = -4826131.800000001
$ws.Range("H136").Value = 2200
$ws.Range("I136").Value = 1433.3334
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 4300.0002
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = 799.9997999999996
$ws.Range("N136").Value = -23700
$ws.Range("H139").Value = 33334870
$ws.Range("I139").Value = 45455604
$ws.Range("J139").Value = 2850
$ws.Range("K139").Value = 136366812
$ws.Range("L139").Value = 8550
$ws.Range("M139").Value = -136361672
$ws.Range("N139").Value = -18830
$ws.Range("H140").Value = 1508.7812
$ws.Range("I140").Value = 1282.5416
$ws.Range("J140").Value = 2187.5
$ws.Range("K140").Value = 3847.6248
$ws.Range("L140").Value = 6562.5
$ws.Range("M140").Value = 1332.3752
$ws.Range("N140").Value = -16922.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 17000
$ws.Range("J137").Value = 17000
$ws.Range("L137").Value = 17000
$ws.Range("N137").Value = -27200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4652548.5
$ws.Range("I132").Value = 8000785.5
$ws.Range("J132").Value = 2219.4443
$ws.Range("K132").Value = 24002356.5
$ws.Range("L132").Value = 6658.3329
$ws.Range("M132").Value = -23999826.5
$ws.Range("N132").Value = -11718.3329
$ws.Range("H136").Value = 18521850
$ws.Range("I136").Value = 22225744
$ws.Range("J136").Value = 2377.7778
$ws.Range("K136").Value = 66677232
$ws.Range("L136").Value = 7133.3334
$ws.Range("M136").Value = -66674682
$ws.Range("N136").Value = -12233.3334
$ws.Range("H137").Value = 20200
$ws.Range("J137").Value = 20200
$ws.Range("L137").Value = 20200
$ws.Range("N137").Value = -30400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 25971372
$ws.Range("I132").Value = 37501296
$ws.Range("J132").Value = 15119679
$ws.Range("K132").Value = 112503888
$ws.Range("L132").Value = 45359037
$ws.Range("M132").Value = -112501358
$ws.Range("N132").Value = -45364097
$ws.Range("H136").Value = 29017886
$ws.Range("I136").Value = 30158228
$ws.Range("J136").Value = 26317076
$ws.Range("K136").Value = 90474684
$ws.Range("L136").Value = 78951228
$ws.Range("M136").Value = -90472134
$ws.Range("N136").Value = -78956328
